$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    4  = -1
    8  = -2
    10 = -4
    17 = -5
    18 = -2
    19 = -6
    22 = 1
    24 = -5
    26 = -8
    28 = -3
    29 = -6
    30 = -2
    32 = 0
    36 = 2
    38 = 2
    39 = -1
    40 = 1
    44 = -4
    45 = -2
    47 = 2
    50 = -1
    54 = 2
    55 = 1
    57 = -2
    58 = -6
    60 = 0
    64 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
